# Fruta / hortaliza, semanal
# Rows 2-11 hold daily price records for "Vega Modelo de Temuco - Frambuesa".
# The underlying data got reshuffled across rows (row 8 stays put); this
# script rewrites the per-row fields (Fecha, Volumen, Precio minimo/maximo/
# promedio, Origen, Precio $/Kg) to match the new arrangement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg) {
    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $PrecioMin
    $ws.Range("O$Row").Value = $PrecioMax
    $ws.Range("P$Row").Value = $PrecioProm
    $ws.Range("R$Row").Value = $Origen
    $ws.Range("S$Row").Value = $PrecioKg
}

Set-Row 2  44176 20  3000 3000 3000 "Región de O'Higgins"    3000
Set-Row 3  44616 200 3200 3200 3200 "Región de La Araucanía" 3200
Set-Row 4  44551 120 4500 4500 4500 "Región de O'Higgins"    4500
Set-Row 5  44574 200 3000 3000 3000 "Región de La Araucanía" 3000
Set-Row 6  44323 20  3200 3200 3200 "Región de La Araucanía" 3200
Set-Row 7  44175 40  5000 5000 5000 "Provincia de Curicó"    5000
Set-Row 9  44214 50  1800 1800 1800 "Región de La Araucanía" 1800
Set-Row 10 44215 65  2800 2800 2800 "Región de La Araucanía" 2800
Set-Row 11 44567 80  2400 2400 2400 "Región de La Araucanía" 2400
